# Cinepolis Maxico Config.xlsx — "maxico input file changed" (week 06 -> week 05)
#
# This script mutates the already-open workbook ($excel.ActiveWorkbook) so that
# the saved .xlsx matches the target revision:
#   * All "06 Envi(o/í)o Semana 06" path strings on the Settings sheet become
#     "05 Envi(o/í)o Semana 05".
#   * WorkingStartDate / WorkingEndDate on the Constants sheet move back one week.
#   * The Settings/Constants sheet scroll position (topLeftCell) and the
#     workbook window geometry are updated to match the author's view state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Settings sheet: update the file-path values (shared strings content).
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("B2").Value  = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos"
$settings.Range("B6").Value  = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos"
$settings.Range("B8").Value  = "/Planeacion/0.Envios TS/2022/05 Envio Semana 05/MÉXICO/Base de Datos/Consolidado"
$settings.Range("B10").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/VIP"
$settings.Range("B12").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Tradicional"
$settings.Range("B14").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Mi Cine"
$settings.Range("B16").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/VIP"
$settings.Range("B20").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos"
$settings.Range("B23").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos/Exportadas"
$settings.Range("B26").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos/Exportadas"
$settings.Range("B29").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos/Exportadas"
$settings.Range("B33").Value = "/Planeacion/0.Envios TS/2022/05 Envío Semana 05/MÉXICO/Base de Datos/Exportadas"

# Scroll the Settings sheet so row 21 is at the top (was row 9).
$settings.Activate() | Out-Null
$settings.Range("B20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 2) Constants sheet: move the working week back by one (dates -7 days).
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

$constants.Range("B24").Value = Get-Date -Year 2022 -Month 1 -Day 24 -Hour 0 -Minute 0 -Second 0
$constants.Range("B25").Value = Get-Date -Year 2022 -Month 1 -Day 30 -Hour 0 -Minute 0 -Second 0

# Scroll the Constants sheet so row 11 is at the top.
$constants.Activate() | Out-Null
$constants.Range("B26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 3) Restore the original active sheet/selection and the workbook window size.
# ---------------------------------------------------------------------------
$settings.Activate() | Out-Null
$settings.Range("B20").Select() | Out-Null

$win = $wb.Windows.Item(1)
$win.Left   = -120
$win.Top    = -120
$win.Width  = 20730
$win.Height = 11160
